# Weekly refresh of the Frambuesa (raspberry) price sheet: a new week's
# record is inserted at row 43 (pushing every following record down by one
# row), so the sheet dimension grows from A1:T148 to A1:T149.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 43:148 down to 44:149 and open up a fresh row 43.
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row with this week's observation.
$ws.Cells.Item(43, 1).Value  = 9
$ws.Cells.Item(43, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(43, 3).Value  = "Metropolitana"
$ws.Cells.Item(43, 4).Value  = 44987
$ws.Cells.Item(43, 5).Value  = 13
$ws.Cells.Item(43, 6).Value  = "Fruta"
$ws.Cells.Item(43, 7).Value  = 100101
$ws.Cells.Item(43, 8).Value  = "Berries"
$ws.Cells.Item(43, 9).Value  = 100101004
$ws.Cells.Item(43, 10).Value = "Frambuesa"
$ws.Cells.Item(43, 11).Value = "Sin especificar"
$ws.Cells.Item(43, 12).Value = "Primera"
$ws.Cells.Item(43, 13).Value = 300
$ws.Cells.Item(43, 14).Value = 7000
$ws.Cells.Item(43, 15).Value = 7000
$ws.Cells.Item(43, 16).Value = 7000
$ws.Cells.Item(43, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(43, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(43, 19).Value = 3500
$ws.Cells.Item(43, 20).Value = 2
